$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Label" column (H) -------------------------------------------
# Copy the formatting of the last existing header cell (G1, bold / bordered /
# centered "style 1") onto the new header cell H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Cells.Item(1, 8).Value = "Label"

# Per-row data: refit D/E (and F where present) values, and the new binary
# diagnosis Label (0 = Control, 1 = MDD) written into column H.
$data = @{
  2  = @{ D = 0.6961575621120364; E = 0.6961575621120364; H = 0 }
  3  = @{ D = 0.4096800711415563; E = 0.4096800711415563; H = 0 }
  4  = @{ D = 0.5907176544974683; E = 0.5907176544974683; H = 0 }
  5  = @{ D = 0.4572414339610804; E = 0.4572414339610804; H = 0 }
  6  = @{ D = 0.5947410258671298; E = 0.5947410258671298; H = 0 }
  7  = @{ D = 0.7021244805004409; E = 0.2978755194995591; H = 1 }
  8  = @{ D = 0.5547843203938406; E = 0.4452156796061594; H = 1 }
  9  = @{ D = 0.4339345447379571; E = 0.566065455262043;  H = 1 }
  10 = @{ D = 0.5187992060299038; E = 0.4812007939700962; H = 1 }
  11 = @{ D = 0.5122682476908;    E = 0.4877317523092;    F = 0.7228833436965942; H = 1 }
  12 = @{ D = 0.7049184343567342; E = 0.7049184343567342; H = 0 }
  13 = @{ D = 0.4783993074252378; E = 0.4783993074252378; H = 0 }
  14 = @{ D = 0.5910940154101465; E = 0.5910940154101465; H = 0 }
  15 = @{ D = 0.4572414339610804; E = 0.4572414339610804; H = 0 }
  16 = @{ D = 0.5954265660426572; E = 0.5954265660426572; H = 0 }
  17 = @{ D = 0.7100399656583939; E = 0.2899600343416061; H = 1 }
  18 = @{ D = 0.5518536802338495; E = 0.4481463197661505; H = 1 }
  19 = @{ D = 0.4503373671123159; E = 0.5496626328876841; H = 1 }
  20 = @{ D = 0.5125855810063282; E = 0.4874144189936718; H = 1 }
  21 = @{ D = 0.5055348271290691; E = 0.4944651728709309; F = 0.7366729974746704; H = 1 }
}

foreach ($r in $data.Keys) {
  $row = $data[$r]
  $ws.Cells.Item($r, 4).Value = $row.D
  $ws.Cells.Item($r, 5).Value = $row.E
  if ($row.ContainsKey("F")) {
    $ws.Cells.Item($r, 6).Value = $row.F
  }
  $ws.Cells.Item($r, 8).Value = $row.H
}
